# Added paths to images in documents
# Sets the image path for the NAWA student document entry (column P, row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Select()
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Student/nawa.jpg"
